$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "parametersDetails" block from I1:M4 to A7:E10 ---
# Header label (single cell)
$ws.Range("I1").Cut($ws.Range("A7"))
# Column headers row
$ws.Range("I2:M2").Cut($ws.Range("A8"))
# Data rows
$ws.Range("I3:M4").Cut($ws.Range("A9"))

# --- Move the "orderdetails"/"measuredValues" block from O1:V3 to A15:H17 ---
# The header text changes from "orderdetails" to "measuredValues", so just clear the old
# header cell and write the new text directly at the destination.
$ws.Range("O1").ClearContents()
$ws.Range("A15").Value = "measuredValues"
# Column headers row
$ws.Range("O2:V2").Cut($ws.Range("A16"))
# Data row
$ws.Range("O3:V3").Cut($ws.Range("A17"))
# R3 (the date cell that moved to D17) leaves a styled-but-empty cell behind; clean it up.
$ws.Range("R3").Clear()

# --- Column widths (best effort; engine quantizes to a 6px grid so exact legacy
# bestFit float widths from real Excel cannot always be reproduced bit-for-bit) ---
$ws.Columns.Item(1).ColumnWidth = 6.0
$ws.Columns.Item(2).ColumnWidth = 10.333333333333334
$ws.Columns.Item(3).ColumnWidth = 6.0
$ws.Columns.Item(6).ColumnWidth = 7.0
$ws.Columns.Item(7).ColumnWidth = 6.833333333333333
$ws.Columns.Item(9).ColumnWidth = 15.166666666666666
$ws.Columns.Item(10).ColumnWidth = 8.666666666666666
$ws.Columns.Item(13).ColumnWidth = 8.666666666666666
$ws.Columns.Item(15).ColumnWidth = 13.666666666666666
$ws.Columns.Item(16).ColumnWidth = 7.0
$ws.Columns.Item(17).ColumnWidth = 11.333333333333334
$ws.Columns.Item(19).ColumnWidth = 3.8333333333333335
$ws.Columns.Item(20).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(21).ColumnWidth = 3.5
$ws.Columns.Item(22).ColumnWidth = 5.0

# --- Sheet view: drop the old scroll position and update the selection ---
[void]$ws.Range("I9").Select()

# --- Page setup: orientation portrait ---
$ws.PageSetup.Orientation = 1
